# BOM-assembly.xlsx edit: add J1, J2, J3 connector rows to BOM, move selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: J2, 2-Pin male header -----------------------------------
$ws.Range("A10").Value = 1
$ws.Range("B10").Value = "2-Pin male header"
$ws.Range("C10").Value = "CONN_1X02"
$ws.Range("D10").Value = "1X02 (Pitch 2.54mm)"
$ws.Range("E10").Value = "J2"
$ws.Range("F10").Value = "2-Pin male header"
$ws.Range("G10").Value = "952-2262-ND"
$ws.Range("H10").Style = "Hyperlink"
$ws.Range("H10").Value = "952-2262-ND"
$ws.Hyperlinks.Add($ws.Range("H10"), "https://www.digikey.com/product-detail/en/harwin-inc/M20-9990246/952-2262-ND/3728226") | Out-Null
$ws.Range("H10").Style = "Hyperlink"
$ws.Range("H10").Value = "https://www.digikey.com/product-detail/en/harwin-inc/M20-9990246/952-2262-ND/3728226"
$ws.Range("I10").Value = " Harwin"
$ws.Range("J10").Value = " M20-9990246"
$ws.Range("K10").Value = 0.11

# --- Row 11: J3, 6-Pin male header -----------------------------------
$ws.Range("A11").Value = 1
$ws.Range("B11").Value = "6-Pin male header"
$ws.Range("C11").Value = "CONN_1X06"
$ws.Range("D11").Value = "1X06 (Pitch 2.54mm)"
$ws.Range("E11").Value = "J3"
$ws.Range("F11").Value = "6-Pin male header"
$ws.Range("G11").Value = "609-3263-ND"
$ws.Range("H11").Style = "Hyperlink"
$ws.Range("H11").Value = "609-3263-ND"
$ws.Hyperlinks.Add($ws.Range("H11"), "https://www.digikey.com/product-detail/en/amphenol-icc-fci/68000-406HLF/609-3263-ND/1878471") | Out-Null
$ws.Range("H11").Style = "Hyperlink"
$ws.Range("H11").Value = "https://www.digikey.com/product-detail/en/amphenol-icc-fci/68000-406HLF/609-3263-ND/1878471"
$ws.Range("I11").Value = " Amphenol"
$ws.Range("J11").Value = "68000-406HLF"
$ws.Range("K11").Value = 0.27

# --- Row 12: J1, 20-pin Female header ---------------------------------
$ws.Range("A12").Value = 1
$ws.Range("B12").Value = "20-pin Female header"
$ws.Range("C12").Value = "CONN_2X10"
$ws.Range("D12").Value = "2X10 (Pitch 2.54mm)"
$ws.Range("E12").Value = "J1"
$ws.Range("F12").Value = "20-pin Female header"
$ws.Range("G12").Value = "S7078-ND"
$ws.Range("H12").Style = "Hyperlink"
$ws.Range("H12").Value = "S7078-ND"
$ws.Hyperlinks.Add($ws.Range("H12"), "https://www.digikey.com/product-detail/en/sullins-connector-solutions/PPTC102LFBN-RC/S7078-ND/810216") | Out-Null
$ws.Range("H12").Style = "Hyperlink"
$ws.Range("H12").Value = "https://www.digikey.com/product-detail/en/sullins-connector-solutions/PPTC102LFBN-RC/S7078-ND/810216"
$ws.Range("I12").Value = "Sullins"
$ws.Range("J12").Value = "PPTC102LFBN-RC"
$ws.Range("K12").Value = 1.3

# --- Misc view state ----------------------------------------------------
$ws.Range("C16").Select() | Out-Null
